$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells already display the raw URL text; read it with .Text (not .Value,
# which this host mis-resolves for Range) so we can wire up a hyperlink whose
# address matches the visible text.
$u4 = $ws.Range("B4").Text
$u5 = $ws.Range("B5").Text

# Insert in the same order the workbook ends up with (B5's link first, then
# B4's) so the <hyperlinks> entries land in that order.
$ws.Hyperlinks.Add($ws.Range("B5"), $u5, "", "", $u5)
$ws.Hyperlinks.Add($ws.Range("B4"), $u4, "", "", $u4)

# Hyperlinks.Add reassigns the "Hipervinculo" cell style, but as a brand new
# style slot; force wrap text back on so the cells land on the same cellXf
# (s="4") that the workbook already used for B3's hyperlink-styled cell.
$ws.Range("B4:B5").WrapText = $true

$ws.Range("B4").Select()
